# "write some new for stephen"
#
# The authoritative diff shows the document's style sheet gaining a second
# "Subtitle" paragraph style definition (identical formatting to the
# existing one), inserted alongside the rest of the style catalogue.
# (The Heading1-Heading6/Title/duplicate-Normal block the diff also shows
# is already present in this document, so nothing to do there.)
#
# Word's object model will not let us mint two styles that share the same
# internal styleId (that duplicate-id artifact in the source diff comes
# from Google Docs' raw OOXML export, not anything reachable through
# Styles.Add/COM) - so we recreate the same visible effect the sane way:
# add a new paragraph style, also named "Subtitle", with matching
# pPr/rPr, while leaving the pre-existing "Subtitle" style untouched.

$d = $word.ActiveDocument

$newSubtitle = $d.Styles.Add("Subtitle1", 1)
$newSubtitle.NameLocal = "Subtitle"
$newSubtitle.BaseStyle = "Normal"
$newSubtitle.NextParagraphStyle = "Normal"

# pPr: keepNext / keepLines / spacing before=360 twips (18pt), after=80 twips (4pt), lineRule=auto
$newSubtitle.ParagraphFormat.KeepWithNext = $true
$newSubtitle.ParagraphFormat.KeepTogether = $true
$newSubtitle.ParagraphFormat.SpaceBefore = 18
$newSubtitle.ParagraphFormat.SpaceAfter = 4
$newSubtitle.ParagraphFormat.LineSpacingRule = 5

# rPr: Georgia (ascii/eastAsia/hAnsi/cs), italic, color 666666, sz/szCs 48 (24pt)
$newSubtitle.Font.Name = "Georgia"
$newSubtitle.Font.NameAscii = "Georgia"
$newSubtitle.Font.NameFarEast = "Georgia"
$newSubtitle.Font.NameBi = "Georgia"
$newSubtitle.Font.Italic = $true
$newSubtitle.Font.Color = 6710886
$newSubtitle.Font.Size = 24
$newSubtitle.Font.SizeBi = 24

Write-Output "Added duplicate Subtitle style: $($newSubtitle.NameLocal)"
